# Adds a new "ViewUserPage" worksheet between "UserPage" and "ResetPage",
# populates it with the view-user detail labels, and adjusts its layout.

$wb = $excel.ActiveWorkbook

$userPage = $wb.Worksheets.Item("UserPage")

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $userPage)
$newSheet.Name = "ViewUserPage"

$newSheet.Range("A1").Value = "Email:"
$newSheet.Range("A2").Value = "Role:"
$newSheet.Range("A3").Value = "Username:"
$newSheet.Range("A4").Value = "Sales Commission Percentage (%):"
$newSheet.Range("A5").Value = "Active"
$newSheet.Range("A6").Value = "Sales Commission Percentage (%):"
$newSheet.Range("A7").Value = "Allowed Contacts: All"

$newSheet.Columns.Item(1).ColumnWidth = 33.17

$null = $newSheet.Range("A7").Select()
